$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D42").Value = "Lime"
$ws.Range("D43").Value = "Shina"
$ws.Range("D44").Value = "Lily"
